$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values on row 3
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 2.6
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 1.95
$ws.Range("X3").Value = 9.5
$ws.Range("Y3").Value = 8.5
$ws.Range("AE3").Value = 15
$ws.Range("AH3").Value = 12
$ws.Range("AJ3").Value = 13
$ws.Range("AO3").Value = 10
$ws.Range("AU3").Value = 8
$ws.Range("AX3").Value = 21

# Add new row 5 with match data
$ws.Range("A5").Value = "0bDmdOdR"
# Copy the Date text from an existing row so it stays plain text (not auto-converted to a date serial)
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("C5").Value = "14:00"
$ws.Range("D5").Value = "LITHUANIA - A LYGA"
$ws.Range("E5").Value = "Zalgiris"
$ws.Range("F5").Value = "Siauliai FA"
$ws.Range("G5").Value = 1.28
$ws.Range("H5").Value = 4.85
$ws.Range("I5").Value = 8.25
$ws.Range("J5").Value = 1.7
$ws.Range("K5").Value = 2.52
$ws.Range("L5").Value = 7.1
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 14.7
$ws.Range("O5").Value = 1.12
$ws.Range("P5").Value = 4.4
$ws.Range("Q5").Value = 1.57
$ws.Range("R5").Value = 2.1
$ws.Range("S5").Value = 1.38
$ws.Range("T5").Value = 2.45
$ws.Range("U5").Value = 1.85
$ws.Range("V5").Value = 1.91
$ws.Range("W5").Value = 7.3
$ws.Range("X5").Value = 6.1
$ws.Range("Y5").Value = 7.3
$ws.Range("Z5").Value = 7.1
$ws.Range("AA5").Value = 8.5
$ws.Range("AB5").Value = 19
$ws.Range("AC5").Value = 15
$ws.Range("AD5").Value = 8.75
$ws.Range("AE5").Value = 16.5
$ws.Range("AF5").Value = 60
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 50
$ws.Range("AJ5").Value = 21
$ws.Range("AK5").Value = 175
$ws.Range("AL5").Value = 75
$ws.Range("AM5").Value = 55
$ws.Range("AN5").Value = 3.2
$ws.Range("AO5").Value = 5.5
$ws.Range("AP5").Value = 14
$ws.Range("AQ5").Value = 13.5
$ws.Range("AR5").Value = 35
$ws.Range("AS5").Value = 175
$ws.Range("AT5").Value = 3.45
$ws.Range("AU5").Value = 8.25
$ws.Range("AV5").Value = 70
$ws.Range("AW5").Value = 9.5
$ws.Range("AX5").Value = 50
$ws.Range("AY5").Value = 45
$ws.Range("AZ5").Value = 350
$ws.Range("BA5").Value = 300
$ws.Range("BB5").Value = 500
